$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 6 de Abril de 2020 a las 01:22'

# Row 4
$ws.Range("B4").Value = 336085
$ws.Range("C4").Value = 24728
$ws.Range("D4").Value = 17245
$ws.Range("E4").Value = 309238
$ws.Range("F4").Value = 8702
$ws.Range("G4").Value = 1151
$ws.Range("H4").Value = 9602

# Row 16
$ws.Range("B16").Value = 15512
$ws.Range("C16").Value = 1600
$ws.Range("D16").Value = 2942
$ws.Range("E16").Value = 12290
$ws.Range("G16").Value = 49
$ws.Range("H16").Value = 280

# Row 52
$ws.Range("A52").Value = 'Argentina'
$ws.Range("B52").Value = 1554
$ws.Range("C52").Value = 103
$ws.Range("D52").Value = 280
$ws.Range("E52").Value = 1228
$ws.Range("F52").Value = 86
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 46

# Row 53
$ws.Range("A53").Value = 'Islandia'
$ws.Range("B53").Value = 1486
$ws.Range("C53").Value = 69
$ws.Range("D53").Value = 428
$ws.Range("E53").Value = 1054
$ws.Range("F53").Value = 11
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 4

# Row 54
$ws.Range("A54").Value = 'Colombia'
$ws.Range("B54").Value = 1485
$ws.Range("C54").Value = 79
$ws.Range("D54").Value = 88
$ws.Range("E54").Value = 1362
$ws.Range("F54").Value = 50
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 35

# Row 87
$ws.Range("B87").Value = 406
$ws.Range("C87").Value = 6
$ws.Range("D87").Value = 104
$ws.Range("E87").Value = 296
$ws.Range("F87").Value = 14
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 6

# Row 117
$ws.Range("A117").Value = 'Mayotte'
$ws.Range("B117").Value = 147
$ws.Range("C117").Value = 13
$ws.Range("D117").Value = 14
$ws.Range("E117").Value = 131
$ws.Range("F117").Value = 3
$ws.Range("H117").Value = 2

# Row 118
$ws.Range("A118").Value = 'Kenia'
$ws.Range("B118").Value = 142
$ws.Range("C118").Value = 16
$ws.Range("D118").Value = 4
$ws.Range("E118").Value = 134
$ws.Range("F118").Value = 2
$ws.Range("H118").Value = 4

# Row 119
$ws.Range("A119").Value = 'Guadalupe'
$ws.Range("B119").Value = 135
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 31
$ws.Range("E119").Value = 97
$ws.Range("F119").Value = 14
$ws.Range("H119").Value = 7

# Row 120
$ws.Range("A120").Value = 'Brunei'
$ws.Range("B120").Value = 135
$ws.Range("D120").Value = 73
$ws.Range("E120").Value = 61
$ws.Range("F120").Value = 3
$ws.Range("H120").Value = 1

# Row 146
$ws.Range("A146").Value = 'Islas Caimanes'
$ws.Range("C146").Value = 4
$ws.Range("E146").Value = 37
$ws.Range("H146").Value = 1

# Row 147
$ws.Range("A147").Value = 'Puerto Rico'
$ws.Range("D147").Value = 1
$ws.Range("E147").Value = 36
$ws.Range("H147").Value = 2

# Row 148
$ws.Range("A148").Value = 'Zambia'
$ws.Range("B148").Value = 39
$ws.Range("D148").Value = 3
$ws.Range("E148").Value = 35
$ws.Range("H148").Value = 1

# Row 149
$ws.Range("A149").Value = 'Bermudas'
$ws.Range("B149").Value = 37
$ws.Range("D149").Value = 14
$ws.Range("E149").Value = 23
$ws.Range("H149").Value = 0

# Row 153
$ws.Range("B153").Value = 29
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 4
$ws.Range("E153").Value = 20
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 5

# Row 159
$ws.Range("A159").Value = 'Haiti'
$ws.Range("C159").Value = 1
$ws.Range("G159").Value = 1

# Row 160
$ws.Range("A160").Value = 'Gabon'
$ws.Range("C160").Value = 0
$ws.Range("G160").Value = 0

# Row 183
$ws.Range("A183").Value = 'Surinam'
$ws.Range("D183").Value = 0
$ws.Range("H183").Value = 1

# Row 184
$ws.Range("A184").Value = 'Mozambique'
$ws.Range("D184").Value = 1
$ws.Range("H184").Value = 0

# Row 192
$ws.Range("A192").Value = 'San Vicente y las Granadinas'

# Row 193
$ws.Range("A193").Value = 'Cabo Verde'
$ws.Range("D193").Value = 0
$ws.Range("H193").Value = 1

# Row 194
$ws.Range("A194").Value = 'Somalia'
$ws.Range("D194").Value = 1
$ws.Range("H194").Value = 0

# Row 196
$ws.Range("A196").Value = 'Botsuana'
$ws.Range("C196").Value = 2
$ws.Range("D196").Value = 0
$ws.Range("H196").Value = 1

# Row 197
$ws.Range("A197").Value = 'Nicaragua'
$ws.Range("C197").Value = 1

# Row 198
$ws.Range("A198").Value = 'San Bartolome'
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 1
$ws.Range("H198").Value = 0

# Row 206
$ws.Range("A206").Value = 'Anguila'

# Row 208
$ws.Range("A208").Value = 'Islas Virgenes Britanicas'

# Row 212
$ws.Range("A212").Value = 'Papua Nueva Guinea'

# Row 213
$ws.Range("A213").Value = 'Timor Oriental'
